# Lab 1 - Results Template: fill in completed experiment data
# (matches commit "Completed data for lab 1")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Motor A PWM test (rows 9-18): updated timing measurements (col T) and
# newly recorded voltage readings (col X)
# ---------------------------------------------------------------------
$ws.Range("T9").Value  = 54.956
$ws.Range("X9").Value  = 6.763

$ws.Range("T10").Value = 11.531
$ws.Range("X10").Value = 6.81

$ws.Range("T11").Value = 7.54
$ws.Range("X11").Value = 6.85

$ws.Range("T12").Value = 6.245
$ws.Range("X12").Value = 6.666

$ws.Range("T13").Value = 4.915
$ws.Range("X13").Value = 6.876

$ws.Range("X14").Value = 6.643
$ws.Range("X15").Value = 6.876
$ws.Range("X16").Value = 6.579
$ws.Range("X17").Value = 6.694
$ws.Range("X18").Value = 6.604

# New "*Using 255 PWM" note above the Motor A voltage column
$ws.Range("W6").Value = "*Using 255 PWM"

# ---------------------------------------------------------------------
# Motor B PWM test (rows 26-36): a new 40%-PWM measurement row was
# inserted ahead of the existing 60/80/100% rows (columns S:U only --
# the Num Washers lifting-test table in columns E:Q is untouched), and
# the Trial/Voltage columns (W:X) were extended with a full run of
# readings (10 trials) one row lower than before, plus a 10th row.
# ---------------------------------------------------------------------

# Row 26 keeps PWM=51(20%) but gets a corrected travel time; the old
# trial-1 data that used to live here moves out (handled below).
$ws.Range("T26").Value = 21.241
$ws.Range("W26").ClearContents()
$ws.Range("X26").ClearContents()

# Row 27 becomes the new 40% PWM row
$ws.Range("S27").Formula = "=255*0.4"
$ws.Range("T27").Value = 7.928
$ws.Range("W27").Value = 1
$ws.Range("X27").Value = 5.989

# Row 28 becomes the 60% PWM row (previously row 27's data)
$ws.Range("S28").Formula = "=255*0.6"
$ws.Range("T28").Value = 5.511
$ws.Range("W28").Value = 2
$ws.Range("X28").Value = 5.663

# Row 29 becomes the 80% PWM row (previously row 28's data)
$ws.Range("S29").Formula = "=255*0.8"
$ws.Range("T29").Value = 4.573
$ws.Range("W29").Value = 3
$ws.Range("X29").Value = 5.77

# Row 30 becomes the 100% PWM row (previously row 29's data)
$ws.Range("S30").Value = 255
$ws.Range("T30").Value = 3.78
$ws.Range("U30").Formula = "=(($C$12/T30)/$C$11)/(2*3.141)*60"
$ws.Range("W30").Value = 4
$ws.Range("X30").Value = 5.823

# Remaining trial rows shift down by one and gain voltage readings
$ws.Range("W31").Value = 5
$ws.Range("X31").Value = 5.952

$ws.Range("W32").Value = 6
$ws.Range("X32").Value = 5.735

$ws.Range("W33").Value = 7
$ws.Range("X33").Value = 5.661

$ws.Range("W34").Value = 8
$ws.Range("X34").Value = 5.57

$ws.Range("W35").Value = 9
$ws.Range("X35").Value = 5.807

$ws.Range("W36").Value = 10
$ws.Range("X36").Value = 5.698

# New "*Using 255 PWM" note above the Motor B voltage column
$ws.Range("W23").Value = "*Using 255 PWM"
